$wb = $excel.ActiveWorkbook

# Rename the "Data_Final" sheet to "Data-Final" to avoid loading errors in Jupyter
$ws = $wb.Worksheets.Item("Data_Final")
$ws.Name = "Data-Final"
